# Applies the coinranking.com 'cryptos' price/volume refresh described by the commit
# "Updated cryptos list on Fri Feb 17 20:26:20 UTC 2023 with GitHub Actions".
#
# Columns: A=Rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# Rows 9/10 and 38/39 also swap two coins' rank position (BinanceUSD<->Polygon,
# WEMIXTOKEN<->Aptos), so B/C (and D/E) are rewritten for those rows too.
#
# Price strings such as "0.9907" / "14.83" read back as pure numbers, so Excel
# would silently reinterpret them as numeric cells on assignment. Prefixing the
# value with a leading apostrophe is the standard COM/UI way to force a literal
# text entry (mirrors typing '0.9907 into the cell) - this matches the source
# workbook, where every Price/Volume cell is stored as text (t="inlineStr").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2: Bitcoin
$ws.Range('D2').Value = '24.775.39'
$ws.Range('E2').Value = '  +0.02%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.704.73'
$ws.Range('E3').Value = '  +0.11%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''0.9907'
$ws.Range('E4').Value = '  -1.38%  '

# Row 5: BNB
$ws.Range('D5').Value = '''313.41'
$ws.Range('E5').Value = '  -1.82%  '

# Row 6: USDC
$ws.Range('D6').Value = '''0.9918'
$ws.Range('E6').Value = '  -1.13%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.3946'
$ws.Range('E7').Value = '  -1.16%  '

# Row 8: Cardano
$ws.Range('D8').Value = '''0.4072'
$ws.Range('E8').Value = '  +1.18%  '

# Row 9: BinanceUSD -> Polygon (rows 9/10 swap)
$ws.Range('B9').Value = 'Polygon'
$ws.Range('C9').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D9').Value = '''1.527'
$ws.Range('E9').Value = '  +7.03%  '

# Row 10: Polygon -> BinanceUSD
$ws.Range('B10').Value = 'BinanceUSD'
$ws.Range('C10').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D10').Value = '''0.9915'
$ws.Range('E10').Value = '  -1.37%  '

# Row 11: OKB
$ws.Range('D11').Value = '''53.88'
$ws.Range('E11').Value = '  +10.61%  '

# Row 12: Dogecoin
$ws.Range('D12').Value = '''0.08769'
$ws.Range('E12').Value = '  -0.82%  '

# Row 13: Polkadot
$ws.Range('D13').Value = '''7.301'
$ws.Range('E13').Value = '  +9.62%  '

# Row 14: Solana
$ws.Range('D14').Value = '''23.30'
$ws.Range('E14').Value = '  -0.56%  '

# Row 15: ShibaInu
$ws.Range('D15').Value = '''0.00001324'
$ws.Range('E15').Value = '  -0.84%  '

# Row 16: Chainlink
$ws.Range('D16').Value = '''7.466'
$ws.Range('E16').Value = '  +2.61%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '1.695.74'
$ws.Range('E17').Value = '  -0.70%  '

# Row 18: Litecoin
$ws.Range('D18').Value = '''100.66'
$ws.Range('E18').Value = '  -2.06%  '

# Row 19: TRON
$ws.Range('D19').Value = '''0.07024'
$ws.Range('E19').Value = '  +1.96%  '

# Row 20: Avalanche
$ws.Range('D20').Value = '''19.51'
$ws.Range('E20').Value = '  -1.40%  '

# Row 21: Uniswap
$ws.Range('D21').Value = '''6.746'
$ws.Range('E21').Value = '  -1.48%  '

# Row 22: Dai
$ws.Range('D22').Value = '''0.9913'
$ws.Range('E22').Value = '  -1.09%  '

# Row 23: Cosmos
$ws.Range('D23').Value = '''14.18'
$ws.Range('E23').Value = '  +1.15%  '

# Row 24: WrappedBTC
$ws.Range('D24').Value = '24.732.73'
$ws.Range('E24').Value = '  -0.07%  '

# Row 25: LidoDAOToken
$ws.Range('D25').Value = '''2.966'
$ws.Range('E25').Value = '  +2.75%  '

# Row 26: Toncoin
$ws.Range('D26').Value = '''2.303'
$ws.Range('E26').Value = '  -0.65%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '''22.36'
$ws.Range('E27').Value = '  -0.02%  '

# Row 28: Monero
$ws.Range('D28').Value = '''158.69'
$ws.Range('E28').Value = '  -1.42%  '

# Row 29: HuobiToken
$ws.Range('D29').Value = '''5.119'
$ws.Range('E29').Value = '  -3.76%  '

# Row 30: BitcoinCash
$ws.Range('D30').Value = '''133.16'
$ws.Range('E30').Value = '  -0.46%  '

# Row 31: Filecoin
$ws.Range('D31').Value = '''7.427'
$ws.Range('E31').Value = '  +26.05%  '

# Row 32: ImmutableX
$ws.Range('D32').Value = '''1.103'
$ws.Range('E32').Value = '  -8.96%  '

# Row 33: WrappedliquidstakedEther2.0
$ws.Range('D33').Value = '1.880.19'
$ws.Range('E33').Value = '  -0.93%  '

# Row 34: Hedera
$ws.Range('D34').Value = '''0.08704'
$ws.Range('E34').Value = '  -4.36%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range('D35').Value = '''7.363'
$ws.Range('E35').Value = '  +19.86%  '

# Row 36: FraxShare
$ws.Range('D36').Value = '''11.18'
$ws.Range('E36').Value = '  +0.65%  '

# Row 37: Algorand
$ws.Range('E37').Value = '  -0.49%  '

# Row 38: WEMIXTOKEN -> Aptos (rows 38/39 swap)
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '''14.83'
$ws.Range('E38').Value = '  -4.50%  '

# Row 39: Aptos -> WEMIXTOKEN
$ws.Range('B39').Value = 'WEMIXTOKEN'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '''1.928'
$ws.Range('E39').Value = '  +3.00%  '

# Row 40: VeChain
$ws.Range('D40').Value = '''0.02781'
$ws.Range('E40').Value = '  +8.67%  '

# Row 41: Stellar
$ws.Range('D41').Value = '''0.08957'
$ws.Range('E41').Value = '  +0.04%  '

# Row 42: TrustWalletToken
$ws.Range('D42').Value = '''1.475'
$ws.Range('E42').Value = '  -0.36%  '

# Row 43: TheSandbox
$ws.Range('D43').Value = '''0.7639'
$ws.Range('E43').Value = '  -0.63%  '

# Row 44: Decentraland
$ws.Range('D44').Value = '''0.7228'
$ws.Range('E44').Value = '  -0.22%  '

# Row 45: EnergySwap
$ws.Range('D45').Value = '''15.41'
$ws.Range('E45').Value = '  -0.67%  '

# Row 46: NEARProtocol
$ws.Range('D46').Value = '''2.459'
$ws.Range('E46').Value = '  -1.29%  '

# Row 47: PancakeSwap
$ws.Range('D47').Value = '''4.145'
$ws.Range('E47').Value = '  -0.67%  '

# Row 48: Frax
$ws.Range('D48').Value = '''0.9909'
$ws.Range('E48').Value = '  -1.18%  '

# Row 49: Quant
$ws.Range('D49').Value = '''140.79'
$ws.Range('E49').Value = '  -1.98%  '

# Row 50: Flow
$ws.Range('D50').Value = '''1.311'
$ws.Range('E50').Value = '  +12.33%  '

# Row 51: Cronos
$ws.Range('D51').Value = '''0.08029'
$ws.Range('E51').Value = '  +0.75%  '
